# Rename the "1PNOTIFY" host/service value to "1PNOTIFYV1" throughout the
# "Notify" worksheet (NotifyTestData.xlsx). This mirrors the commit
# "Modified 1PNotify service into 1PNotifyV1": every cell that previously
# held the literal text "1PNOTIFY" should now read "1PNOTIFYV1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notify")

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $usedRange.Cells.Item($r, $c)
        if ($cell.Value2 -eq "1PNOTIFY") {
            $cell.Value2 = "1PNOTIFYV1"
        }
    }
}

# Move the active selection to L5 (matches the saved view state in the diff).
$ws.Range("L5").Select()
